$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$commentText = $ws.Range("F1").Comment.Text()
$ws.Range("F1").Comment.Delete()

$ws.Columns.Item(6).Insert()

$ws.Cells.Item(1, 6).Value = "Product Group"
$ws.Cells.Item(1, 6).HorizontalAlignment = -4108
$ws.Cells.Item(1, 6).Font.Bold = $true
$ws.Cells.Item(1, 6).Font.Size = 12

$ws.Range("G1").AddComment($commentText) | Out-Null

$ws.Range("A1:Q1").Borders.Item(9).Weight = -4138
$ws.Range("G1:Q1").Font.Size = 12
$ws.Range("G1:Q1").Font.Bold = $true
